$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (Excel would otherwise auto-convert numeric-looking text into real numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.795.53'
$ws.Range("E2").Value = '  +7.07%  '
$ws.Range("D3").Value = '1.993.22'
$ws.Range("E3").Value = '  +5.74%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '252.84'
$ws.Range("E5").Value = '  +3.04%  '
$ws.Range("D6").Value = '0.701'
$ws.Range("E6").Value = '  +2.30%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '48.48'
$ws.Range("E8").Value = '  +13.87%  '
$ws.Range("E9").Value = '  +8.93%  '
$ws.Range("D10").Value = '59.34'
$ws.Range("E10").Value = '  +7.59%  '
$ws.Range("D11").Value = '0.0773'
$ws.Range("E11").Value = '  +4.41%  '
$ws.Range("E12").Value = '  +2.44%  '
$ws.Range("D13").Value = '15.71'
$ws.Range("E13").Value = '  +14.19%  '
$ws.Range("D14").Value = '0.846'
$ws.Range("E14").Value = '  +9.04%  '
$ws.Range("D15").Value = '2.258.64'
$ws.Range("E15").Value = '  +4.55%  '
$ws.Range("D16").Value = '5.23'
$ws.Range("E16").Value = '  +4.88%  '
$ws.Range("D17").Value = '1.965.71'
$ws.Range("E17").Value = '  +3.82%  '
$ws.Range("D18").Value = '38.013.95'
$ws.Range("E18").Value = '  +7.69%  '
$ws.Range("D19").Value = '75.81'
$ws.Range("E19").Value = '  +3.55%  '
$ws.Range("D20").Value = '0.0₃0868'
$ws.Range("E20").Value = '  +5.55%  '
$ws.Range("D21").Value = '13.86'
$ws.Range("E21").Value = '  +8.48%  '
$ws.Range("D22").Value = '255.33'
$ws.Range("E22").Value = '  +4.67%  '
$ws.Range("D23").Value = '5.28'
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("D25").Value = '2.51'
$ws.Range("E25").Value = '  -6.41%  '
$ws.Range("D26").Value = '170.27'
$ws.Range("E26").Value = '  +1.79%  '
$ws.Range("D27").Value = '2.16'
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("D28").Value = '8.99'
$ws.Range("E28").Value = '  +5.54%  '
$ws.Range("D29").Value = '19.18'
$ws.Range("E29").Value = '  +5.13%  '
$ws.Range("E30").Value = '  +1.98%  '
$ws.Range("D31").Value = '4.62'
$ws.Range("E31").Value = '  +7.98%  '
$ws.Range("D32").Value = '0.0621'
$ws.Range("E32").Value = '  +4.80%  '
$ws.Range("D33").Value = '0.0922'
$ws.Range("E33").Value = '  +27.53%  '
$ws.Range("D34").Value = '4.38'
$ws.Range("E34").Value = '  +5.10%  '
$ws.Range("B35").Value = 'Gas'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D35").Value = '20.38'
$ws.Range("E35").Value = '  +54.03%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.91'
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '0.907'
$ws.Range("E38").Value = '  +6.85%  '
$ws.Range("D39").Value = '1.47'
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("D40").Value = '2.08'
$ws.Range("E40").Value = '  +7.60%  '
$ws.Range("D41").Value = '105.26'
$ws.Range("E41").Value = '  +7.65%  '
$ws.Range("D42").Value = '0.0230'
$ws.Range("E42").Value = '  +4.16%  '
$ws.Range("D43").Value = '17.79'
$ws.Range("E43").Value = '  +3.99%  '
$ws.Range("D44").Value = '2.90'
$ws.Range("E44").Value = '  +20.50%  '
$ws.Range("E45").Value = '  +6.04%  '
$ws.Range("D46").Value = '1.373.27'
$ws.Range("E46").Value = '  +3.27%  '
$ws.Range("E47").Value = '  +3.55%  '
$ws.Range("D48").Value = '0.0853'
$ws.Range("E48").Value = '  +5.55%  '
$ws.Range("D49").Value = '2.85'
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("D50").Value = '3.96'
$ws.Range("E50").Value = '  +19.55%  '
$ws.Range("D51").Value = '6.52'
$ws.Range("E51").Value = '  +4.05%  '

# Restore the default cell style for the Price column so no stray
# number-format style is left behind on the saved workbook.
$ws.Range("D2:D51").Style = "Normal"
